$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 "Employment volumes": next-period release date refreshed
$ws.Range("D2").Value = "Apr 2022 - Mar 2023 (15/08/23)"

# Row 3 "Employment by occupation": now coded with SOC2020, so the
# latest-period figure (previously blocked on the SOC coding issue) is
# available and the next-period date is refreshed. The custom date number
# format that used to decorate the latest-period cell is no longer needed.
$ws.Range("C3").Value = "Jan 2022 - Dec 2022 (07/23)"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "Apr 2022 - Mar 2023 (15/08/23)"

# Row 4 "Employment by industry": next-period release date refreshed
$ws.Range("D4").Value = "Apr 2022 - Mar 2023 (15/08/23)"
